$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
try {
    Write-Host "Design.Name:" $d.Name
} catch {
    Write-Host "get err:" $_.Exception.Message
}
try {
    $d.Name = "Office Theme"
    Write-Host "set design name ok, now:" $d.Name
} catch {
    Write-Host "set err:" $_.Exception.Message
}

$m = $p.SlideMaster
try {
    Write-Host "Master.Name:" $m.Name
    $m.Name = "Office Theme"
    Write-Host "Master.Name after:" $m.Name
} catch {
    Write-Host "master name err:" $_.Exception.Message
}
